$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new contact record (Shah, Imran) was added as the second data row of the
# credentials table - insert a fresh row right below the header so every
# existing row shifts down by one.
$ws.Rows.Item(2).Insert()

# Fill in the new row with the new account's data (columns A-G).
$ws.Cells.Item(2, 1).Value = "Shah"
$ws.Cells.Item(2, 2).Value = "Imran"
$ws.Cells.Item(2, 3).Value = "mrshahimranshovon@gmail.com"
$ws.Cells.Item(2, 4).Value = "twg5ZmvBqsixAfr"
$ws.Cells.Item(2, 5).Value = "185.24.233.182:4006"
$ws.Cells.Item(2, 6).Value = "8GbKtEpRUr29jbg6"
$ws.Cells.Item(2, 7).Value = "TMwprA4NyqSKxc6V"

# Widen the PROXY:PORT / PROXY_USER columns so the longer values are readable.
$ws.Columns.Item(5).ColumnWidth = 24.1666
$ws.Columns.Item(6).ColumnWidth = 31.8333

# Give the header row an explicit height (matches the saved workbook) and
# leave the active cell parked on the new row's password column.
$ws.Rows.Item(1).RowHeight = 14.4
$ws.Range("D2").Select() | Out-Null
